# Add a new row (row 7) of match data to the active sheet, mirroring
# the existing flat layout of Id/Date/Time/League/Home/Away + odds columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 'O6ibYFEq'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '07/11/2024'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = '14:00'
$ws.Range("D7").Value = 'SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE'
$ws.Range("E7").Value = 'Al Orubah'
$ws.Range("F7").Value = 'Al Ittihad'
$ws.Range("G7").Value = 6.5
$ws.Range("H7").Value = 4.75
$ws.Range("I7").Value = 1.4
$ws.Range("J7").Value = 6
$ws.Range("K7").Value = 2.5
$ws.Range("L7").Value = 1.83
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 10.5
$ws.Range("O7").Value = 1.17
$ws.Range("P7").Value = 4.5
$ws.Range("Q7").Value = 1.57
$ws.Range("R7").Value = 2.35
$ws.Range("S7").Value = 1.29
$ws.Range("T7").Value = 3.5
$ws.Range("U7").Value = 1.8
$ws.Range("V7").Value = 1.91
$ws.Range("W7").Value = 21
$ws.Range("X7").Value = 41
$ws.Range("Y7").Value = 21
$ws.Range("Z7").Value = 67
$ws.Range("AA7").Value = 41
$ws.Range("AB7").Value = 41
$ws.Range("AC7").Value = 15
$ws.Range("AD7").Value = 9.5
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 51
$ws.Range("AG7").Value = 500
$ws.Range("AH7").Value = 8.5
$ws.Range("AI7").Value = 7.5
$ws.Range("AJ7").Value = 8.5
$ws.Range("AK7").Value = 10
$ws.Range("AL7").Value = 11
$ws.Range("AM7").Value = 23
$ws.Range("AN7").Value = 8
$ws.Range("AO7").Value = 34
$ws.Range("AP7").Value = 34
$ws.Range("AQ7").Value = 101
$ws.Range("AR7").Value = 126
$ws.Range("AS7").Value = 400
$ws.Range("AT7").Value = 3.5
$ws.Range("AU7").Value = 8.5
$ws.Range("AV7").Value = 51
$ws.Range("AW7").Value = 3.6
$ws.Range("AX7").Value = 7
$ws.Range("AY7").Value = 17
$ws.Range("AZ7").Value = 17
$ws.Range("BA7").Value = 41
$ws.Range("BB7").Value = 101
$ws.Range("BC7").Value = 81
$ws.Range("BD7").Value = 81
